# Update the "想去人数" (want-to-go count) figures in column F for both the
# "展览" and "全部类型" worksheets, matching the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 393
$ws1.Range("F3").Value = 119
$ws1.Range("F4").Value = 1635
$ws1.Range("F5").Value = 17
$ws1.Range("F6").Value = 23
$ws1.Range("F9").Value = 0

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F6").Value = 0
$ws4.Range("F7").Value = 0
$ws4.Range("F8").Value = 145
$ws4.Range("F10").Value = 507
